$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# D2: switch the "Name/Firma..." stammdaten field-spec from the old ":pflicht"
# separator syntax to the new "|pflicht" syntax.
$ws.Range("D2").Value = "Name/Firma:text|pflicht;Anschrift:text|pflicht;Steuernummer:text|pflicht;AMA-Betriebsnummer:text;SVS-Versicherungsnummer:text;Bankverbindung(IBAN/BIC):text;Familienstand:select(ledig,verh.,geschieden,verwitwet)|pflicht;Kinder:number;Vollmacht-und-DSGVO-Einwilligung:checkbox|pflicht"

# D3: "Erwerbsart" checkbox field-spec switches to the new "|" separator syntax
# and becomes optional instead of pflicht, and drops the space after the comma.
$ws.Range("D3").Value = "Erwerbsart:checkbox(Vollerwerb,Nebenerwerb)|optional"

# D5: the small "Alle Daten hochgeladen?/Upload" field-spec also moves to the
# new "|pflicht" syntax (other cells still sharing the old text, e.g. D7:D10
# and D17, are left untouched).
$ws.Range("D5").Value = "Alle Daten hochgeladen?:info;Upload:checkbox|pflicht"

# Move the active selection from D2 to D5, matching the saved sheet view.
$ws.Range("D5").Select()
